# Edit material list: add "Дата" (Date) column, update several rows,
# and remove the last row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove row 5 entirely (shifts nothing else, it's the last row) ---
$ws.Rows.Item(5).Delete()

# --- Add new header cell I1 "Дата", copying the header formatting from H1 ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "Дата"

# --- Row 2 values ---
$ws.Range("B2").Value = "Sjsnisb"
$ws.Range("C2").Value = "кг"
$ws.Range("D2").Value = "'5"
$ws.Range("E2").Value = "суммы"
$ws.Range("F2").Value = "'12"
$ws.Range("G2").Value = "'60"
$ws.Range("H2").Value = "effrfwwpkp"
$ws.Range("I2").Value = "22.04.2021"

# --- Row 3 values ---
$ws.Range("B3").Value = "fgrgrg"
$ws.Range("C3").Value = "м^2"
$ws.Range("D3").Value = "'12"
$ws.Range("E3").Value = "суммы"
$ws.Range("F3").Value = "'1"
$ws.Range("G3").Value = "'12"
$ws.Range("H3").Value = "effrfwwpkp"
$ws.Range("I3").Value = "22.04.2021"

# --- Row 4 values ---
$ws.Range("B4").Value = "'122122112"
$ws.Range("C4").Value = "м"
$ws.Range("D4").Value = "'1"
$ws.Range("E4").Value = "суммы"
$ws.Range("F4").Value = "'80"
$ws.Range("G4").Value = "'80"
$ws.Range("H4").Value = "effrfwwpkp"
$ws.Range("I4").Value = "22.04.2021"
